# stock_orderpoint_impex_matrix / bad_wh_name.xlsx
# Commit: "[WIP] Add default_code column, update matching to default_code
#          Update column order Update tests Update formatting of cells"
#
# The fixture gains a new first data column ("Code Article" / default_code)
# to the left of the existing "Article" column. All existing columns shift
# one place to the right, and the new A column is populated with the
# product's internal reference ("E-COM11") on the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every existing column one place to the right (B1->C1, A2->B2, ...)
# and make room for the new leading "Code Article" column.
$ws.Columns("A:A").Insert()

# New header cell for the inserted column.
$ws.Range("A2").Value2 = "Code Article"

# New data cell: the product's internal reference/default_code.
$ws.Range("A3").Value2 = "E-COM11"

# Restore the (now stale) selection to match the refreshed sheet state.
[void]$ws.Range("E10").Select()
